$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate a Range whose text exactly equals $searchText (literal,
# case-sensitive match over the whole document).
# ---------------------------------------------------------------------------
function Find-ExactRange([string]$searchText) {
    $r = $d.Content.Duplicate
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    return $r
}

# ---------------------------------------------------------------------------
# Helper: replace the text of a Range with $newText while keeping that run
# isolated from its (identically formatted) neighbouring runs. The engine
# auto-coalesces adjacent runs that end up with identical formatting, which
# would otherwise merge our freshly edited run into whatever sits next to
# it. Temporarily toggling Bold on, editing the text, then toggling Bold
# back off forces the run to stay split out as its own <w:r>.
# ---------------------------------------------------------------------------
function Set-RangeTextIsolated($range, [string]$newText) {
    $start = $range.Start
    $range.Bold = 1
    $editRange = $d.Range($start, $range.End)
    $editRange.Text = $newText
    $resultRange = $d.Range($start, $start + $newText.Length)
    $resultRange.Bold = 0
    return $resultRange
}

# ---------------------------------------------------------------------------
# 1) Description meta tag: reword the intro sentence, drop the trailing
#    period from the run that follows it (". " -> " ").
# ---------------------------------------------------------------------------
$introRange = Find-ExactRange("is our first article in a new series of Game Maker, where we will be introducing you to it")
Set-RangeTextIsolated $introRange "article introduces the room, and how to add elements or game pieces to it." | Out-Null

$periodRange = Find-ExactRange(". ")
Set-RangeTextIsolated $periodRange " " | Out-Null

# ---------------------------------------------------------------------------
# 2) Revised meta tag: replace "Thursday 30th, 2025" (spanning 3 runs, one
#    of them superscript "th") with a single run "November 19, 2025".
# ---------------------------------------------------------------------------
$revisedRange = Find-ExactRange("Thursday 30th, 2025")
Set-RangeTextIsolated $revisedRange "November 19, 2025" | Out-Null

# ---------------------------------------------------------------------------
# 3) Url meta tag: point at the new "4_The_Room" article instead of the old
#    "1_Introduction_to_Game_Maker" one.
# ---------------------------------------------------------------------------
$urlRange = Find-ExactRange("Enlightenment/Articles/2025/4_Game_Maker/1_Introduction_to_Game_Maker/1_Introduction_to_GameMaker.html")
Set-RangeTextIsolated $urlRange "Enlightenment/Articles/2025/4_Game_Maker/4_The_Room/4_The_Room.html" | Out-Null

Write-Output "Done."
